# Update two-digit multiplication problems throughout the document.
# Each "old" string is unique in the document, so a straightforward
# Find/Replace (MatchWholeWord) per pair is safe and deterministic.

$d = $word.ActiveDocument

$replacements = @(
    @("66×29=", "68×43="),
    @("58×96=", "76×51="),
    @("46×75=", "57×40="),
    @("43×84=", "86×15="),
    @("17×37=", "22×42="),
    @("48×53=", "38×33="),
    @("65×74=", "77×14="),
    @("27×19=", "30×11="),
    @("43×38=", "55×93="),
    @("36×42=", "74×19="),
    @("35×18=", "25×36="),
    @("86×31=", "35×70="),
    @("54×63=", "23×45="),
    @("12×93=", "33×73="),
    @("66×50=", "61×19="),
    @("63×67=", "53×75="),
    @("62×66=", "26×12="),
    @("88×41=", "82×75="),
    @("83×49=", "98×38="),
    @("13×34=", "32×82="),
    @("60×63=", "88×45="),
    @("17×77=", "93×18="),
    @("13×48=", "74×63="),
    @("24×48=", "89×85="),
    @("17×51=", "51×54=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
